# regen sval data to filter save games
# Overwrite the B2:G13 numeric data block with the regenerated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.006876353814593728,0.004309184025731883,16.98373111632243,0.4998867070740569,1,17.49480336123681),
    @(3.182878228561681,1.65323645889881,0.7127328510149897,6.48142807727062,0,12.0302756157461),
    @(0.1554434735375247,0.004309184025731883,0.7127328510149897,0.4998867070740569,0,1.372372215652303),
    @(0.1554434735375247,0.3375848360084654,3.082599426703578,0.4998867070740569,0,4.075514443323626),
    @(1.505614041169197,1.65323645889881,0.7127328510149897,0.4998867070740569,1,4.371470058157054),
    @(3.182878228561681,1.65323645889881,0.1529057820181812,0.4998867070740569,1,5.488907176552729),
    @(3.182878228561681,1.65323645889881,0.7127328510149897,0.4998867070740569,0,6.048734245549538),
    @(0.7287194209349384,1.65323645889881,0.7127328510149897,0.4998867070740569,0,3.594575437922795),
    @(0.7287194209349384,1.65323645889881,0.1529057820181812,0.4998867070740569,0,3.034748368925986),
    @(3.182878228561681,1.65323645889881,0.1529057820181812,0.4998867070740569,0,5.488907176552729),
    @(0.7287194209349384,1.65323645889881,0.1529057820181812,0.4998867070740569,0,3.034748368925986),
    @(1.505614041169197,1.65323645889881,3.082599426703578,0.4998867070740569,0,6.741336633845642)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $c = $j + 2
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}
